# Add new algorithm timetable entry (row) into Sheet1 and snapshot the
# previous Sheet1 contents into Sheet2 (used as the "before" data for the
# new graph-coloring algorithm work).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Copy Sheet1's current rows (1-5, cols A-L) into Sheet2 first,
#        before Sheet1 gets its new row inserted. ---
for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 12; $c++) {
        $srcCell = $ws1.Cells.Item($r, $c)
        $dstCell = $ws2.Cells.Item($r, $c)
        $dstCell.Value = $srcCell.Value2
    }
}

# --- 2. Insert a new row at row 2 on Sheet1 (shifts existing rows 2-5 to 3-6) ---
$ws1.Rows.Item(2).Insert()

# --- 3. Populate the newly inserted row 2 with the new course entry ---
$ws1.Cells.Item(2, 1).Value = "elec"
$ws1.Cells.Item(2, 2).Value = 221
$ws1.Cells.Item(2, 3).Value = 202
$ws1.Cells.Item(2, 4).Value = "L2A"
$ws1.Cells.Item(2, 5).Value = "T2A"
$ws1.Cells.Item(2, 6).Value = "D2B"
$ws1.Cells.Item(2, 7).Value = "THURS"
$ws1.Cells.Item(2, 8).Value = 1900
$ws1.Cells.Item(2, 9).Value = 2100

# --- 4. Fix up selections to match the final, saved state.
#        Sheet2's selection is set first, then Sheet1 is reselected last so
#        Sheet1 stays the active/visible tab (as in the original workbook). ---
[void]$ws2.Range("A1:L5").Select()
[void]$ws1.Select()
[void]$ws1.Range("G2").Select()
